# "Analiza 2 kazalcev, poprava tabel"
# The "Operating revenue (Turnover)" line (row 3) is removed from the
# financial table; every row below it shifts up by one. This is the
# same effect as selecting row 3 and deleting the entire row in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Operating revenue (Turnover)" row (row 3) entirely -
# cells below shift up, row heights travel with their row, and the
# shared string for the removed label is dropped once unused.
$ws.Rows(3).Delete()

# Leave the selection on the (new) row 3, matching the post-edit
# selection left by the "select row, delete row" workflow.
$ws.Range("A3:XFD3").Select()
